$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "41.523.94",
# "0.0790"); force text format first so Excel does not silently
# coerce them to floating point numbers and lose exact formatting.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.523.94'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.468.22'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.74'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.82'
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  +2.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.24'
$ws.Range('E10').Value = '  -4.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.848.73'
$ws.Range('E13').Value = '  -0.75%  '
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.88'
$ws.Range('E15').Value = '  +2.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.470.78'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.772'
$ws.Range('E17').Value = '  -2.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.517.35'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('E19').Value = '  +1.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('E21').Value = '  +3.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.05'
$ws.Range('E22').Value = '  -2.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.66'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('E26').Value = '  -0.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.64'
$ws.Range('E27').Value = '  +1.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  -1.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.32'
$ws.Range('E30').Value = '  -3.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.83'
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.57'
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -4.74%  '
$ws.Range('E36').Value = '  -7.23%  '
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  -6.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.21'
$ws.Range('E40').Value = '  -11.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.01'
$ws.Range('E41').Value = '  -5.38%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.945.29'
$ws.Range('E43').Value = '  -2.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0283'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.44'
$ws.Range('E45').Value = '  -7.44%  '
$ws.Range('E46').Value = '  -3.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.02'
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.706.77'
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '96.82'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '67.08'
$ws.Range('E50').Value = '  -4.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.15'
$ws.Range('E51').Value = '  +2.04%  '
